# First working version of batch simulation.
$wb = $excel.ActiveWorkbook

# --- Rename the original sheet ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "cell types"

# --- Lowercase the existing row labels (rows that do not move) ---
$ws1.Range("A2").Value = "species"
$ws1.Range("A3").Value = "proliferation rate"
$ws1.Range("A4").Value = "abundance"
$ws1.Range("A5").Value = "metabolism"

# --- Insert a new "food to move" row before the existing "food to divide" row (old row 6) ---
$ws1.Rows.Item(6).Insert() | Out-Null

$ws1.Range("A6").Value = "food to move"
$ws1.Range("B6").Formula = "=`$B`$5*4"
$ws1.Range("C6").Formula = "=`$C`$5*4"
$ws1.Range("D6").Formula = "=`$D`$5*4"

# --- Relabel "food to divide" (now shifted to row 7) ---
$ws1.Range("A7").Value = "food to divide"

# --- Remove the old "food to survive" (*10) row, originally row 7, now row 8 after the insert ---
$ws1.Rows.Item(8).Delete() | Out-Null

# --- Relabel/update "division recovery time" row (now row 8), values 20 -> 5 ---
$ws1.Range("A8").Value = "division recovery time"
$ws1.Range("B8").Value = 5
$ws1.Range("C8").Value = 5
$ws1.Range("D8").Value = 5

# --- Relabel remaining rows ---
$ws1.Range("A9").Value = "food to survive"
$ws1.Range("A10").Value = "endurance"

# --- Clear tab-selected state on sheet1, set the selection to B10 ---
$ws1.Range("B10").Select() | Out-Null

# --- Add the new sheet after "cell types" ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "190824-00"

$ws2.Range("B1").Value = "S001"
$ws2.Range("C1").Value = "S002"
$ws2.Range("D1").Value = "S003"

$ws2.Range("A2").Value = "width"
$ws2.Range("B2:D2").Value = 500

$ws2.Range("A3").Value = "height"
$ws2.Range("B3:D3").Value = 500

$ws2.Range("A4").Value = "maxIter"
$ws2.Range("B4:D4").Value = 10

$ws2.Range("A5").Value = "seeds"
$ws2.Range("B5:D5").Value = 20

$ws2.Range("A6").Value = "foodFile"
$ws2.Range("B6:D6").Value = "foodMaps-04.png"

$ws2.Range("A7").Value = "mapFile"
$ws2.Range("B7:D7").Value = "foodMaps-00.png"

$ws2.Range("A8").Value = "cellTypeNames"
$ws2.Range("B8:D8").Value = "C001, C002, C003"

$ws2.Range("A9").Value = "mixRatios"
$ws2.Range("B9:D9").Value = "1, 1, 1"

$ws2.Range("A10").Value = "outputSize"
$ws2.Range("B10:D10").Value = "2160, 2160"

$ws2.Range("A11").Value = "timeWarpFactor"
$ws2.Range("B11:D11").Value = 2

# --- Formatting: center align everything in B:D, plus text format for rows 6-10 ---
$ws2.Range("B1:D11").HorizontalAlignment = -4108  # xlCenter
$ws2.Range("B6:D10").NumberFormat = "@"

# --- Column widths (target: col A ~14.66 chars bestFit, cols B:D ~20.83 chars) ---
$ws2.Columns.Item(1).ColumnWidth = 13.83
$ws2.Range("B1:D1").ColumnWidth = 20.0

# --- Selection on new sheet ---
$ws2.Range("B4:D4").Select() | Out-Null
